# Update countries & provincias Spain
# Daily refresh of the COVID country dataset on sheet "Pais":
#  - A handful of neighbouring countries swapped ranking order (their case
#    counts crossed over), so their name/data pairs are exchanged between
#    two adjacent rows.
#  - A number of countries simply received updated totals for the day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Estados Unidos (row 4) : updated totals ---
$ws.Range("B4").Value = 6969917
$ws.Range("C4").Value = 2514
$ws.Range("D4").Value = 4223987
$ws.Range("E4").Value = 2542079
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 203851

# --- India (row 5) : updated totals ---
$ws.Range("B5").Value = 5412565
$ws.Range("C5").Value = 14335
$ws.Range("D5").Value = 4310645
$ws.Range("E5").Value = 1015027
$ws.Range("G5").Value = 119
$ws.Range("H5").Value = 86893

# --- Alemania (row 25) : updated totals ---
$ws.Range("B25").Value = 272342
$ws.Range("C25").Value = 34
$ws.Range("E25").Value = 19376

# --- Catar (row 32) : updated totals ---
$ws.Range("B32").Value = 123376
$ws.Range("C32").Value = 230
$ws.Range("D32").Value = 120303
$ws.Range("E32").Value = 2863
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 210

# --- Kuwait (row 40) : updated totals ---
$ws.Range("B40").Value = 99434
$ws.Range("C40").Value = 385
$ws.Range("D40").Value = 90168
$ws.Range("E40").Value = 8682
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 584

# --- Oman / Paises Bajos (rows 41-42) swap order with new data ---
$ws.Range("A41").Value = "Paises Bajos"
$ws.Range("B41").Value = 93778
$ws.Range("C41").Value = 1844
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 6279

$ws.Range("A42").Value = "Oman"
$ws.Range("B42").Value = 93475
$ws.Range("C42").Value = 1722
$ws.Range("D42").Value = 85418
$ws.Range("E42").Value = 7211
$ws.Range("G42").Value = 28
$ws.Range("H42").Value = 846

# --- Afganistan / Azerbaiyan (rows 67-68) swap order with new data ---
$ws.Range("A67").Value = "Azerbaiyan"
$ws.Range("B67").Value = 39188
$ws.Range("C67").Value = 146
$ws.Range("D67").Value = 36755
$ws.Range("E67").Value = 1858
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 575

$ws.Range("A68").Value = "Afganistan"
$ws.Range("B68").Value = 39044
$ws.Range("C68").Value = 125
$ws.Range("D68").Value = 32576
$ws.Range("E68").Value = 5027
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 1441

# --- Dinamarca (row 81) : updated totals ---
$ws.Range("B81").Value = 22905
$ws.Range("C81").Value = 469
$ws.Range("D81").Value = 17514
$ws.Range("E81").Value = 4753
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 638

# --- Senegal (row 90) : updated totals ---
$ws.Range("B90").Value = 14714
$ws.Range("C90").Value = 26
$ws.Range("D90").Value = 11260
$ws.Range("E90").Value = 3152

# --- Birmania (row 115) : updated totals ---
$ws.Range("E115").Value = 3986
$ws.Range("G115").Value = 8
$ws.Range("H115").Value = 89

# --- Timor Oriental / Santa Lucia (rows 204-205) swap order, data unchanged ---
$ws.Range("A204").Value = "Santa Lucia"
$ws.Range("A205").Value = "Timor Oriental"

# --- Islas Malvinas / Montserrat (rows 214-215) swap order with new data ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
